$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2's user record is swapped from "jhonsonbaby / Jhonson Babay" to
# "vladimir / Vladimir Jhonson" (uname, firstName and lastName change).
$ws.Range("A2").Value = "vladimir"
$ws.Range("D2").Value = "Vladimir"
$ws.Range("E2").Value = "Jhonson"

# The sheet had the same record duplicated on row 3 (with its own copy of
# the emailAddress hyperlink) - drop that duplicate row entirely.
# Hyperlinks.Delete() on this host clears the *whole* worksheet's hyperlink
# collection (it isn't scoped to the calling range), so remove them all,
# delete the now-unwanted row, then restore the single hyperlink that
# should remain, on F2.
$ws.Range("F2").Hyperlinks.Delete()
$ws.Rows("3:3").Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:jhonson@yopmail.con")

# Re-adding the hyperlink stamps F2 with the built-in "Hyperlink" cell
# style; put it back to Normal and drop that now-unused style definition so
# every cell in the used range ends up sharing one plain, explicit style.
$ws.Range("F2").Style = "Normal"
$wb.Styles("Hyperlink").Delete()

# Touch the whole used range so each cell picks up an explicit style index
# (mirrors the workbook's new, single duplicated cellXfs entry) instead of
# relying on the default (no "s" attribute).
$ws.Range("A1:G2").IndentLevel = 0

# Update the selection to cover the (now smaller) used range.
$ws.Range("A1:G2").Select()
